{"js": "// Replace each two-digit-multiplication equation in the worksheet table\n// with the newly generated equation/answer pair. Each cell contains a\n// single unique run of text, so searching the whole document body for\n// each old equation string and replacing it in place is sufficient and\n// unambiguous.\n\nconst mapping = [\n  [\"91\u00d783=7553\", \"23\u00d764=1472\"],\n  [\"76\u00d714=1064\", \"83\u00d753=4399\"],\n  [\"75\u00d713=975\", \"69\u00d728=1932\"],\n  [\"92\u00d776=6992\", \"26\u00d721=546\"],\n  [\"93\u00d753=4929\", \"82\u00d742=3444\"],\n  [\"34\u00d759=2006\", \"78\u00d746=3588\"],\n  [\"96\u00d744=4224\", \"67\u00d747=3149\"],\n  [\"19\u00d754=1026\", \"94\u00d759=5546\"],\n  [\"47\u00d713=611\", \"59\u00d731=1829\"],\n  [\"85\u00d713=1105\", \"21\u00d771=1491\"],\n  [\"75\u00d725=1875\", \"52\u00d765=3380\"],\n  [\"26\u00d762=1612\", \"74\u00d755=4070\"],\n  [\"71\u00d792=6532\", \"32\u00d742=1344\"],\n  [\"14\u00d786=1204\", \"58\u00d786=4988\"],\n  [\"82\u00d748=3936\", \"49\u00d716=784\"],\n  [\"90\u00d714=1260\", \"76\u00d716=1216\"],\n  [\"95\u00d790=8550\", \"61\u00d761=3721\"],\n  [\"39\u00d743=1677\", \"48\u00d796=4608\"],\n  [\"56\u00d714=784\", \"40\u00d791=3640\"],\n  [\"63\u00d747=2961\", \"21\u00d780=1680\"],\n  [\"60\u00d739=2340\", \"80\u00d747=3760\"],\n  [\"42\u00d724=1008\", \"35\u00d777=2695\"],\n  [\"78\u00d783=6474\", \"78\u00d751=3978\"],\n  [\"97\u00d770=6790\", \"77\u00d765=5005\"],\n  [\"90\u00d728=2520\", \"99\u00d747=4653\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of mapping) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-multiplication equation in the worksheet table\n# with the newly generated equation/answer pair. Each cell contains a\n# single unique run of text, so a straightforward Find/Replace per pair\n# is sufficient and unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"91\u00d783=7553\", \"23\u00d764=1472\"),\n    @(\"76\u00d714=1064\", \"83\u00d753=4399\"),\n    @(\"75\u00d713=975\",  \"69\u00d728=1932\"),\n    @(\"92\u00d776=6992\", \"26\u00d721=546\"),\n    @(\"93\u00d753=4929\", \"82\u00d742=3444\"),\n    @(\"34\u00d759=2006\", \"78\u00d746=3588\"),\n    @(\"96\u00d744=4224\", \"67\u00d747=3149\"),\n    @(\"19\u00d754=1026\", \"94\u00d759=5546\"),\n    @(\"47\u00d713=611\",  \"59\u00d731=1829\"),\n    @(\"85\u00d713=1105\", \"21\u00d771=1491\"),\n    @(\"75\u00d725=1875\", \"52\u00d765=3380\"),\n    @(\"26\u00d762=1612\", \"74\u00d755=4070\"),\n    @(\"71\u00d792=6532\", \"32\u00d742=1344\"),\n    @(\"14\u00d786=1204\", \"58\u00d786=4988\"),\n    @(\"82\u00d748=3936\", \"49\u00d716=784\"),\n    @(\"90\u00d714=1260\", \"76\u00d716=1216\"),\n    @(\"95\u00d790=8550\", \"61\u00d761=3721\"),\n    @(\"39\u00d743=1677\", \"48\u00d796=4608\"),\n    @(\"56\u00d714=784\",  \"40\u00d791=3640\"),\n    @(\"63\u00d747=2961\", \"21\u00d780=1680\"),\n    @(\"60\u00d739=2340\", \"80\u00d747=3760\"),\n    @(\"42\u00d724=1008\", \"35\u00d777=2695\"),\n    @(\"78\u00d783=6474\", \"78\u00d751=3978\"),\n    @(\"97\u00d770=6790\", \"77\u00d765=5005\"),\n    @(\"90\u00d728=2520\", \"99\u00d747=4653\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
